$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-03 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-04 Monday", 2)

$d.Content.Find.Execute("821÷9=91, 2", $true, $false, $false, $false, $false, $true, 1, $false, "388÷7=55, 3", 2)
$d.Content.Find.Execute("863÷4=215, 3", $true, $false, $false, $false, $false, $true, 1, $false, "265÷8=33, 1", 2)
$d.Content.Find.Execute("828÷9=92, 0", $true, $false, $false, $false, $false, $true, 1, $false, "383÷9=42, 5", 2)
$d.Content.Find.Execute("424÷7=60, 4", $true, $false, $false, $false, $false, $true, 1, $false, "162÷8=20, 2", 2)
$d.Content.Find.Execute("240÷2=120, 0", $true, $false, $false, $false, $false, $true, 1, $false, "733÷5=146, 3", 2)

$d.Content.Find.Execute("925÷2=462, 1", $true, $false, $false, $false, $false, $true, 1, $false, "153÷6=25, 3", 2)
$d.Content.Find.Execute("241÷3=80, 1", $true, $false, $false, $false, $false, $true, 1, $false, "322÷9=35, 7", 2)
$d.Content.Find.Execute("613÷2=306, 1", $true, $false, $false, $false, $false, $true, 1, $false, "233÷4=58, 1", 2)
$d.Content.Find.Execute("272÷2=136, 0", $true, $false, $false, $false, $false, $true, 1, $false, "809÷5=161, 4", 2)
$d.Content.Find.Execute("924÷6=154, 0", $true, $false, $false, $false, $false, $true, 1, $false, "400÷9=44, 4", 2)

$d.Content.Find.Execute("567÷5=113, 2", $true, $false, $false, $false, $false, $true, 1, $false, "558÷5=111, 3", 2)
$d.Content.Find.Execute("833÷7=119, 0", $true, $false, $false, $false, $false, $true, 1, $false, "936÷9=104, 0", 2)
$d.Content.Find.Execute("105÷2=52, 1", $true, $false, $false, $false, $false, $true, 1, $false, "198÷2=99, 0", 2)
$d.Content.Find.Execute("704÷9=78, 2", $true, $false, $false, $false, $false, $true, 1, $false, "228÷7=32, 4", 2)
$d.Content.Find.Execute("978÷3=326, 0", $true, $false, $false, $false, $false, $true, 1, $false, "483÷7=69, 0", 2)

$d.Content.Find.Execute("688÷8=86, 0", $true, $false, $false, $false, $false, $true, 1, $false, "525÷3=175, 0", 2)
$d.Content.Find.Execute("287÷6=47, 5", $true, $false, $false, $false, $false, $true, 1, $false, "306÷9=34, 0", 2)
$d.Content.Find.Execute("259÷2=129, 1", $true, $false, $false, $false, $false, $true, 1, $false, "216÷9=24, 0", 2)
$d.Content.Find.Execute("949÷9=105, 4", $true, $false, $false, $false, $false, $true, 1, $false, "865÷9=96, 1", 2)
$d.Content.Find.Execute("570÷4=142, 2", $true, $false, $false, $false, $false, $true, 1, $false, "464÷5=92, 4", 2)

$d.Content.Find.Execute("990÷8=123, 6", $true, $false, $false, $false, $false, $true, 1, $false, "427÷6=71, 1", 2)
$d.Content.Find.Execute("441÷6=73, 3", $true, $false, $false, $false, $false, $true, 1, $false, "254÷5=50, 4", 2)
$d.Content.Find.Execute("461÷8=57, 5", $true, $false, $false, $false, $false, $true, 1, $false, "977÷6=162, 5", 2)
$d.Content.Find.Execute("587÷9=65, 2", $true, $false, $false, $false, $false, $true, 1, $false, "485÷3=161, 2", 2)
$d.Content.Find.Execute("446÷7=63, 5", $true, $false, $false, $false, $false, $true, 1, $false, "479÷4=119, 3", 2)
